# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labelled columns AC:AE -----------------------
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting already used by the rest of the header row
# (bold font, thin box border, centered horizontal / top vertical alignment).
$header = $ws.Range("AC1:AE1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# --- Data rows (rows 2-46): team's season record on every player row ------
$ws.Range("AC2:AC46").Value = 65
$ws.Range("AD2:AD46").Value = 97
$ws.Range("AE2:AE46").Value = 0
